$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 286, shifting existing rows 286..303 down to 287..304.
$ws.Rows.Item(286).Insert()

# Populate the new row 286 with the new record (matching style of surrounding rows).
$ws.Cells.Item(286, 1).Value = 5
$ws.Cells.Item(286, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(286, 3).Value = "Maule"
$ws.Cells.Item(286, 4).Value = 44714
$ws.Cells.Item(286, 5).Value = 7
$ws.Cells.Item(286, 6).Value = 100112003
$ws.Cells.Item(286, 7).Value = "Ajo"
$ws.Cells.Item(286, 8).Value = "Chino"
$ws.Cells.Item(286, 9).Value = "Primera"
$ws.Cells.Item(286, 10).Value = 300
$ws.Cells.Item(286, 11).Value = 18000
$ws.Cells.Item(286, 12).Value = 18000
$ws.Cells.Item(286, 13).Value = 18000
$ws.Cells.Item(286, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(286, 15).Value = "China"
$ws.Cells.Item(286, 16).Value = 1800
$ws.Cells.Item(286, 17).Value = 10
$ws.Cells.Item(286, 18).Value = "Hortaliza"
